# Update Leve profit calculation figures across multiple job sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 446.875
$ws.Range("I9").Value = 342.30768
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 342.30768
$ws.Range("L9").Value = 900
$ws.Range("M9").Value = -173.30768
$ws.Range("N9").Value = -1238

$ws.Range("H58").Value = 22321.18
$ws.Range("I58").Value = 376.25
$ws.Range("J58").Value = 24229.436
$ws.Range("K58").Value = 1128.75
$ws.Range("L58").Value = 72688.308
$ws.Range("M58").Value = -978.75
$ws.Range("N58").Value = -72988.308

$ws.Range("H74").Value = 4124
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 4137.778
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 4137.778
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -6009.778

$ws.Range("H77").Value = 4124
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 4137.778
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 20688.89
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -30048.89

$ws.Range("H103").Value = 562.0833
$ws.Range("I103").Value = 508
$ws.Range("J103").Value = 600.7143
$ws.Range("K103").Value = 1524
$ws.Range("L103").Value = 1802.1429
$ws.Range("M103").Value = -938
$ws.Range("N103").Value = -2974.1429

$ws.Range("H115").Value = 2250
$ws.Range("I115").Value = 500
$ws.Range("K115").Value = 1500
$ws.Range("M115").Value = 67

$ws.Range("H138").Value = 4558.5845
$ws.Range("I138").Value = 3077.2222
$ws.Range("J138").Value = 4796.6606
$ws.Range("K138").Value = 9231.6666
$ws.Range("L138").Value = 14389.9818
$ws.Range("M138").Value = -4091.6666
$ws.Range("N138").Value = -24669.9818

$ws.Range("H141").Value = 351577.78
$ws.Range("I141").Value = 1149.25
$ws.Range("K141").Value = 3447.75
$ws.Range("M141").Value = 1732.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 27455.666
$ws.Range("I25").Value = 1174.5
$ws.Range("J25").Value = 80018
$ws.Range("K25").Value = 1174.5
$ws.Range("L25").Value = 80018
$ws.Range("M25").Value = -772.5
$ws.Range("N25").Value = -80822

$ws.Range("H33").Value = 16142.714
$ws.Range("I33").Value = 8999.5
$ws.Range("K33").Value = 8999.5
$ws.Range("M33").Value = -8670.5

$ws.Range("H102").Value = 3666.6667
$ws.Range("I102").Value = 2875
$ws.Range("K102").Value = 2875
$ws.Range("M102").Value = -1253

$ws.Range("H132").Value = 12822432
$ws.Range("I132").Value = 17858142
$ws.Range("J132").Value = 4263.273
$ws.Range("K132").Value = 53574426
$ws.Range("L132").Value = 12789.819
$ws.Range("M132").Value = -53571896
$ws.Range("N132").Value = -17849.819

$ws.Range("H137").Value = 29642.857
$ws.Range("J137").Value = 29642.857
$ws.Range("L137").Value = 29642.857
$ws.Range("N137").Value = -39842.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4566.7295
$ws.Range("I31").Value = 3072.8262
$ws.Range("J31").Value = 7021
$ws.Range("K31").Value = 3072.8262
$ws.Range("L31").Value = 7021
$ws.Range("M31").Value = -2777.8262
$ws.Range("N31").Value = -7611

$ws.Range("H34").Value = 4566.7295
$ws.Range("I34").Value = 3072.8262
$ws.Range("J34").Value = 7021
$ws.Range("K34").Value = 3072.8262
$ws.Range("L34").Value = 7021
$ws.Range("M34").Value = -2870.8262
$ws.Range("N34").Value = -7425

$ws.Range("H74").Value = 17579.182
$ws.Range("J74").Value = 17579.182
$ws.Range("L74").Value = 17579.182
$ws.Range("N74").Value = -19327.182

$ws.Range("H77").Value = 17579.182
$ws.Range("J77").Value = 17579.182
$ws.Range("L77").Value = 52737.546
$ws.Range("N77").Value = -61473.546

$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 10000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -12996

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1967
$ws.Range("J75").Value = 2921.25
$ws.Range("L75").Value = 8763.75
$ws.Range("N75").Value = -10759.75

$ws.Range("H78").Value = 1967
$ws.Range("J78").Value = 2921.25
$ws.Range("L78").Value = 26291.25
$ws.Range("N78").Value = -36275.25

$ws.Range("H87").Value = 11820
$ws.Range("J87").Value = 15816.667
$ws.Range("L87").Value = 47450.001
$ws.Range("N87").Value = -49946.001

$ws.Range("H90").Value = 11820
$ws.Range("J90").Value = 15816.667
$ws.Range("L90").Value = 142350.003
$ws.Range("N90").Value = -154830.003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 74178.36
$ws.Range("I102").Value = 2044.2222
$ws.Range("J102").Value = 204019.8
$ws.Range("K102").Value = 2044.2222
$ws.Range("L102").Value = 204019.8
$ws.Range("M102").Value = -422.2221999999999
$ws.Range("N102").Value = -207263.8

$ws.Range("H122").Value = 4520
$ws.Range("I122").Value = 2678.5715
$ws.Range("J122").Value = 7098
$ws.Range("K122").Value = 8035.7145
$ws.Range("L122").Value = 21294
$ws.Range("M122").Value = -5585.7145
$ws.Range("N122").Value = -26194

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4310
$ws.Range("I61").Value = 1922
$ws.Range("J61").Value = 6300
$ws.Range("K61").Value = 1922
$ws.Range("L61").Value = 6300
$ws.Range("M61").Value = -1720
$ws.Range("N61").Value = -6704

$ws.Range("H113").Value = 4310
$ws.Range("I113").Value = 1922
$ws.Range("J113").Value = 6300
$ws.Range("K113").Value = 1922
$ws.Range("L113").Value = 6300
$ws.Range("M113").Value = 248
$ws.Range("N113").Value = -10640

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 33333.332
$ws.Range("J125").Value = 33333.332
$ws.Range("L125").Value = 33333.332
$ws.Range("N125").Value = -43173.332

Write-Host "Updated Leve profit rows across ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets."
